# Generate Report for Handoff
#
# Replaces the old (png-based) sample rows with four markdown-based rows
# (calleeMd1.md, calleeMd2.md, callerMd1.md, callerMd2.md) across the
# Overview / zh-cn / de-de sheets, and appends a 4th data row to each
# (rows 2-5 instead of 2-4).

$wb = $excel.ActiveWorkbook
$newline = [char]10

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 2; File = "calleeMd1.md" },
    @{ Row = 3; File = "calleeMd2.md" },
    @{ Row = 4; File = "callerMd1.md" },
    @{ Row = 5; File = "callerMd2.md" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $ws1.Range("A$row").Value = $r.File
    $ws1.Range("B$row").Value = "Ready for handoff"
    $ws1.Range("C$row").Value = "Ready for handoff"
    $ws1.Range("D$row").Value = "2016-03-30 10:54:45"
}

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null

# Re-apply the link look (font colour gets reset to the theme hyperlink
# colour by Hyperlinks.Add) and the date/time number format, in one pass
# at the end so the stylesheet doesn't end up with one-off duplicates.
$ws1.Range("A2:A5").Font.Underline = 2
$ws1.Range("A2:A5").Font.Color = 15570276
$ws1.Range("D2:D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhRows = @(
    @{ Row = 2; File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"; K = "e2e\callerMd2.md,${newline}e2e\callerMd1.md"; I = $null },
    @{ Row = 3; File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"; K = "e2e\callerMd1.md"; I = $null },
    @{ Row = 4; File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"; K = $null; I = "e2e\calleeMd1.md,${newline}e2e\calleeMd2.md" },
    @{ Row = 5; File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"; K = $null; I = "e2e\calleeMd1.md" }
)

foreach ($r in $zhRows) {
    $row = $r.Row

    $ws2.Range("A$row").Value = $r.File
    $ws2.Range("B$row").Value = ".md"
    $ws2.Range("C$row").Value = "Ready for handoff"
    $ws2.Range("D$row").Value = $r.Xlf
    $ws2.Range("E$row").Value = "2016-03-30 10:54:34"
    $ws2.Range("H$row").Value = "0001-01-01 00:00:00"

    if ($r.I) {
        $ws2.Range("I$row").Value = $r.I
    }

    $ws2.Range("J$row").Value = "Include"

    if ($r.K) {
        $ws2.Range("K$row").Value = $r.K
    }
}

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf") | Out-Null

$ws2.Range("A2:A5").Font.Underline = 2
$ws2.Range("A2:A5").Font.Color = 15570276
$ws2.Range("D2:D5").Font.Underline = 2
$ws2.Range("D2:D5").Font.Color = 15570276
$ws2.Range("E2:E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "de-de"
# =======================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$deRows = @(
    @{ Row = 2; File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"; K = "e2e\callerMd2.md,${newline}e2e\callerMd1.md"; I = $null },
    @{ Row = 3; File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"; K = "e2e\callerMd1.md"; I = $null },
    @{ Row = 4; File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"; K = $null; I = "e2e\calleeMd1.md,${newline}e2e\calleeMd2.md" },
    @{ Row = 5; File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"; K = $null; I = "e2e\calleeMd1.md" }
)

foreach ($r in $deRows) {
    $row = $r.Row

    $ws3.Range("A$row").Value = $r.File
    $ws3.Range("B$row").Value = ".md"
    $ws3.Range("C$row").Value = "Ready for handoff"
    $ws3.Range("D$row").Value = $r.Xlf
    $ws3.Range("E$row").Value = "2016-03-30 10:54:45"
    $ws3.Range("H$row").Value = "0001-01-01 00:00:00"

    if ($r.I) {
        $ws3.Range("I$row").Value = $r.I
    }

    $ws3.Range("J$row").Value = "Include"

    if ($r.K) {
        $ws3.Range("K$row").Value = $r.K
    }
}

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf") | Out-Null

$ws3.Range("A2:A5").Font.Underline = 2
$ws3.Range("A2:A5").Font.Color = 15570276
$ws3.Range("D2:D5").Font.Underline = 2
$ws3.Range("D2:D5").Font.Color = 15570276
$ws3.Range("E2:E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Report regenerated."
